$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-08 Monday" "2024-04-09 Tuesday"

Replace-Text "88÷3=" "16÷7="
Replace-Text "68÷6=" "88÷4="
Replace-Text "44÷6=" "81÷9="
Replace-Text "89÷5=" "79÷9="
Replace-Text "83÷2=" "39÷7="
Replace-Text "94÷4=" "91÷3="
Replace-Text "88÷5=" "66÷6="
Replace-Text "13÷7=" "34÷6="
Replace-Text "99÷5=" "44÷6="
Replace-Text "61÷4=" "39÷3="
Replace-Text "60÷8=" "19÷2="
Replace-Text "56÷8=" "59÷8="
Replace-Text "81÷5=" "40÷3="
Replace-Text "95÷9=" "82÷6="
Replace-Text "93÷7=" "93÷5="
Replace-Text "86÷9=" "51÷8="
Replace-Text "58÷8=" "57÷7="
Replace-Text "56÷9=" "54÷6="
Replace-Text "18÷6=" "21÷7="
Replace-Text "16÷5=" "23÷6="
Replace-Text "59÷7=" "79÷4="
Replace-Text "34÷3=" "52÷9="
Replace-Text "25÷3=" "56÷6="
Replace-Text "54÷7=" "69÷8="
Replace-Text "60÷2=" "60÷7="
